# Remove the generated "Ver no Jupiter..." line, the trailing
# copyright/footer line, and the blank paragraph that separated them
# from the body text. They used to follow the
# "Devido as caracteristicas praticas..." paragraph and precede the
# trailing blank paragraph / page break.

$d = $word.ActiveDocument

$anchorText = "Devido às características práticas da disciplina, não será oferecida recuperação"
$footerText = "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

$anchorPara = $null
$footerPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq $anchorText) {
        $anchorPara = $p
    }
    if ($txt -eq $footerText) {
        $footerPara = $p
    }
}

if ($anchorPara -ne $null -and $footerPara -ne $null) {
    $startPara = $anchorPara.Next()
    $deleteRange = $d.Range($startPara.Range.Start, $footerPara.Range.End)
    $deleteRange.Delete()
}
